# Append three new rows (132-134) to Sheet1, each a copy of the last
# existing data row (131) with the date in column A advanced by one day.
# All other columns (B:J) keep the same values as row 131.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 131
$newDates = @(45688, 45689, 45690)

# Values for columns B..J, copied verbatim from row 131.
$rowValues = @(116.4121952, 0.00170247, 0.008850780000000001, 0.06933635, 12792.90181321, 465.80531254, 0.24, 1.7904431, 485.38834923)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $destRow = $lastRow + 1 + $i

    # Copy formatting (style) from the template row so the new date cell
    # keeps the same style (border/alignment/date number format) as the
    # rest of column A, without introducing new style entries.
    $ws.Range("A$lastRow" + ":J$lastRow").Copy()
    $ws.Range("A$destRow" + ":J$destRow").PasteSpecial(-4122)

    # Write the actual values for the new row.
    $ws.Range("A$destRow").Value = $newDates[$i]
    $ws.Range("B$destRow").Value = $rowValues[0]
    $ws.Range("C$destRow").Value = $rowValues[1]
    $ws.Range("D$destRow").Value = $rowValues[2]
    $ws.Range("E$destRow").Value = $rowValues[3]
    $ws.Range("F$destRow").Value = $rowValues[4]
    $ws.Range("G$destRow").Value = $rowValues[5]
    $ws.Range("H$destRow").Value = $rowValues[6]
    $ws.Range("I$destRow").Value = $rowValues[7]
    $ws.Range("J$destRow").Value = $rowValues[8]
}

$excel.CutCopyMode = 0
